$wb = $excel.ActiveWorkbook

# ----- Sheet: ALC -----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H88").Value = 1506.9375
$ws.Range("I88").Value = 841
$ws.Range("J88").Value = 1906.5
$ws.Range("K88").Value = 841
$ws.Range("L88").Value = 1906.5
$ws.Range("M88").Value = -435
$ws.Range("N88").Value = -2718.5
$ws.Range("H91").Value = 1506.9375
$ws.Range("I91").Value = 841
$ws.Range("J91").Value = 1906.5
$ws.Range("K91").Value = 841
$ws.Range("L91").Value = 1906.5
$ws.Range("M91").Value = 563
$ws.Range("N91").Value = -4714.5
$ws.Range("H100").Value = 2905
$ws.Range("I100").Value = 2905
$ws.Range("K100").Value = 2905
$ws.Range("M100").Value = -2364
$ws.Range("H132").Value = 1859.56
$ws.Range("I132").Value = 1825.6086
$ws.Range("J132").Value = 2250
$ws.Range("K132").Value = 5476.825800000001
$ws.Range("L132").Value = 6750
$ws.Range("M132").Value = -2946.825800000001
$ws.Range("N132").Value = -11810
$ws.Range("H135").Value = 1004
$ws.Range("I135").Value = 635.6
$ws.Range("K135").Value = 5720.400000000001
$ws.Range("M135").Value = -3185.400000000001
$ws.Range("H137").Value = 2040.7273
$ws.Range("I137").Value = 986.8
$ws.Range("K137").Value = 2960.4
$ws.Range("M137").Value = -410.3999999999996

# ----- Sheet: ARM -----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 0
$ws.Range("I61").Value = 0
$ws.Range("K61").Value = 0
$ws.Range("M61").ClearContents()
$ws.Range("H122").Value = 3617.182
$ws.Range("I122").Value = 1957.8
$ws.Range("K122").Value = 5873.4
$ws.Range("M122").Value = -3423.4
$ws.Range("H132").Value = 4446.2856
$ws.Range("I132").Value = 1781
$ws.Range("K132").Value = 5343
$ws.Range("M132").Value = -2813
$ws.Range("H136").Value = 0
$ws.Range("I136").Value = 0
$ws.Range("K136").Value = 0
$ws.Range("M136").ClearContents()

# ----- Sheet: BSM -----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 1139.5
$ws.Range("I94").Value = 1139.5
$ws.Range("K94").Value = 1139.5
$ws.Range("M94").Value = -688.5

# ----- Sheet: CRP -----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3880.75
$ws.Range("I31").Value = 1684.7142
$ws.Range("K31").Value = 1684.7142
$ws.Range("M31").Value = -1389.7142
$ws.Range("H34").Value = 3880.75
$ws.Range("I34").Value = 1684.7142
$ws.Range("K34").Value = 1684.7142
$ws.Range("M34").Value = -1482.7142
$ws.Range("H39").Value = 0
$ws.Range("I39").Value = 0
$ws.Range("K39").Value = 0
$ws.Range("M39").ClearContents()
$ws.Range("H49").Value = 0
$ws.Range("I49").Value = 0
$ws.Range("K49").Value = 0
$ws.Range("M49").ClearContents()
$ws.Range("H99").Value = 11071.871
$ws.Range("I99").Value = 8661.637000000001
$ws.Range("J99").Value = 12397.5
$ws.Range("K99").Value = 8661.637000000001
$ws.Range("L99").Value = 12397.5
$ws.Range("M99").Value = -7163.637000000001
$ws.Range("N99").Value = -15393.5
$ws.Range("H122").Value = 5443.6665
$ws.Range("I122").Value = 5608.8125
$ws.Range("J122").Value = 5113.375
$ws.Range("K122").Value = 16826.4375
$ws.Range("L122").Value = 15340.125
$ws.Range("M122").Value = -14376.4375
$ws.Range("N122").Value = -20240.125
$ws.Range("H126").Value = 11071.871
$ws.Range("I126").Value = 8661.637000000001
$ws.Range("J126").Value = 12397.5
$ws.Range("K126").Value = 25984.911
$ws.Range("L126").Value = 37192.5
$ws.Range("M126").Value = -23514.911
$ws.Range("N126").Value = -42132.5
$ws.Range("H132").Value = 4464.7334
$ws.Range("I132").Value = 2802.3333
$ws.Range("K132").Value = 8406.999899999999
$ws.Range("M132").Value = -5876.999899999999

# ----- Sheet: CUL -----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H62").Value = 2349.75
$ws.Range("J62").Value = 3300
$ws.Range("L62").Value = 9900
$ws.Range("N62").Value = -11272
$ws.Range("H65").Value = 2349.75
$ws.Range("J65").Value = 3300
$ws.Range("L65").Value = 29700
$ws.Range("N65").Value = -36564
$ws.Range("H68").Value = 1802.5
$ws.Range("I68").Value = 2699
$ws.Range("J68").Value = 1578.375
$ws.Range("K68").Value = 8097
$ws.Range("L68").Value = 4735.125
$ws.Range("M68").Value = -7286
$ws.Range("N68").Value = -6357.125
$ws.Range("J69").Value = 5000
$ws.Range("L69").Value = 15000
$ws.Range("N69").Value = -16622
$ws.Range("H70").Value = 894
$ws.Range("I70").Value = 894
$ws.Range("J70").Value = 0
$ws.Range("K70").Value = 2682
$ws.Range("L70").Value = 0
$ws.Range("M70").Value = -2367
$ws.Range("H71").Value = 1802.5
$ws.Range("I71").Value = 2699
$ws.Range("J71").Value = 1578.375
$ws.Range("K71").Value = 24291
$ws.Range("L71").Value = 14205.375
$ws.Range("M71").Value = -20235
$ws.Range("N71").Value = -22317.375
$ws.Range("J72").Value = 5000
$ws.Range("L72").Value = 45000
$ws.Range("N72").Value = -53112
$ws.Range("H73").Value = 894
$ws.Range("I73").Value = 894
$ws.Range("J73").Value = 0
$ws.Range("K73").Value = 2682
$ws.Range("L73").Value = 0
$ws.Range("M73").Value = -1590
$ws.Range("H75").Value = 266.5
$ws.Range("I75").Value = 266.5
$ws.Range("K75").Value = 799.5
$ws.Range("M75").Value = 198.5
$ws.Range("H78").Value = 266.5
$ws.Range("I78").Value = 266.5
$ws.Range("K78").Value = 2398.5
$ws.Range("M78").Value = 2593.5
$ws.Range("H80").Value = 5999
$ws.Range("J80").Value = 5999
$ws.Range("L80").Value = 17997
$ws.Range("N80").Value = -19869
$ws.Range("H83").Value = 5999
$ws.Range("J83").Value = 5999
$ws.Range("L83").Value = 53991
$ws.Range("N83").Value = -63351
$ws.Range("H122").Value = 1893
$ws.Range("J122").Value = 0
$ws.Range("L122").Value = 0
$ws.Range("N122").ClearContents()

# ----- Sheet: GSM -----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 6353.6
$ws.Range("I97").Value = 6742
$ws.Range("K97").Value = 6742
$ws.Range("M97").Value = -6246
$ws.Range("H102").Value = 1540.091
$ws.Range("I102").Value = 522.4706
$ws.Range("K102").Value = 522.4706
$ws.Range("M102").Value = 1099.5294
$ws.Range("H113").Value = 2262.577
$ws.Range("J113").Value = 3881.0833
$ws.Range("L113").Value = 3881.0833
$ws.Range("N113").Value = -8221.0833
$ws.Range("H122").Value = 36992.867
$ws.Range("J122").Value = 87096.5
$ws.Range("L122").Value = 261289.5
$ws.Range("N122").Value = -266189.5
$ws.Range("H132").Value = 4054.75
$ws.Range("I132").Value = 2887.5
$ws.Range("J132").Value = 5222
$ws.Range("K132").Value = 8662.5
$ws.Range("L132").Value = 15666
$ws.Range("M132").Value = -6132.5
$ws.Range("N132").Value = -20726

# ----- Sheet: LTW -----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H100").Value = 5662.3335
$ws.Range("I100").Value = 993.5
$ws.Range("J100").Value = 15000
$ws.Range("K100").Value = 993.5
$ws.Range("L100").Value = 15000
$ws.Range("M100").Value = -452.5
$ws.Range("N100").Value = -16082
$ws.Range("H136").Value = 3556.5

# ----- Sheet: WVR -----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H55").Value = 30000
$ws.Range("I55").Value = 30000
$ws.Range("K55").Value = 30000
$ws.Range("M55").Value = -29723
$ws.Range("H107").Value = 899.8125
$ws.Range("I107").Value = 500.66666
$ws.Range("K107").Value = 1501.99998
$ws.Range("M107").Value = 418.0000199999999
